# Update "want-to-go" count values (column F) on both the "展览" sheet
# and the "全部类型" sheet, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 14652
$ws1.Range("F5").Value = 17879
$ws1.Range("F24").Value = 7377
$ws1.Range("F28").Value = 1186
$ws1.Range("F35").Value = 233

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 14652
$ws4.Range("F5").Value = 17879
$ws4.Range("F25").Value = 7377
$ws4.Range("F29").Value = 1186
$ws4.Range("F37").Value = 233
